# Femacal de La Calera - Sandia: add a new weekly price-report row at the
# top of the data block (row 721), pushing all existing rows down by one.
# This mirrors a new week's data being inserted at the front of the table,
# with the rest of the historical rows shifting down (and the table
# growing by one row overall, from R758 to R759).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 721:758 down to 722:759, leaving a blank row 721.
$ws.Rows.Item(721).Insert()

# Populate the new row 721 with the new week's record.
$ws.Cells.Item(721, 1).Value  = 3
$ws.Cells.Item(721, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(721, 3).Value  = "Coquimbo"
$ws.Cells.Item(721, 4).Value  = 45267
$ws.Cells.Item(721, 5).Value  = 5
$ws.Cells.Item(721, 6).Value  = 100112028
$ws.Cells.Item(721, 7).Value  = "Sandia"
$ws.Cells.Item(721, 8).Value  = "Sin especificar"
$ws.Cells.Item(721, 9).Value  = "Primera"
$ws.Cells.Item(721, 10).Value = 350
$ws.Cells.Item(721, 11).Value = 800
$ws.Cells.Item(721, 12).Value = 800
$ws.Cells.Item(721, 13).Value = 800
$ws.Cells.Item(721, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(721, 15).Value = "Perú"
$ws.Cells.Item(721, 16).Value = 800
$ws.Cells.Item(721, 17).Value = 1
$ws.Cells.Item(721, 18).Value = "Hortaliza"
